$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "(-1)For incorrectly adding products to the customer who already exists."

$ws.Range("F20").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 3
